$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 7
$ws.Range("C2").Value = 0.5119816233366274
$ws.Range("D2").Value = 0.2013802160866826
$ws.Range("F2").Value = [double]"2.781931710304319e-22"

# Row 3
$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 0.5027128607879965
$ws.Range("D3").Value = 0.1278603665557776
$ws.Range("F3").Value = [double]"5.497000879062199e-28"

# Row 4
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 0.3942212794016117
$ws.Range("D4").Value = 0.1649984488824481
$ws.Range("F4").Value = [double]"1.032691392870033e-21"
